# Simulated Wild Card round and logged it
# Updates "R" (road) row Target Depth Data on both OFF and DEF sheets
# for the Patriots, reflecting the results of the Wild Card round game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 307
$wsOff.Range("C3").Value = 218
$wsOff.Range("D3").Value = 146
$wsOff.Range("E3").Value = 61
$wsOff.Range("F3").Value = 6
$wsOff.Range("G3").Value = 8

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 460
$wsDef.Range("C3").Value = 325
$wsDef.Range("D3").Value = 122
$wsDef.Range("E3").Value = 48
